$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tarifas bancárias")

# --- BB (Banco do Brasil) bank-fee entries ---
$ws.Range("A3").Value = "GRUPO XYZ"
$ws.Range("B3").Value = "BB"
$ws.Range("C3").Value = 1234
$ws.Range("D3").Value = "4567-8"
$ws.Range("E3").Value = "09/01/2025"
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = "TARIFA - BB"

$ws.Range("A4").Value = "GRUPO XYZ"
$ws.Range("B4").Value = "BB"
$ws.Range("C4").Value = 1234
$ws.Range("D4").Value = "4567-8"
$ws.Range("E4").Value = "09/02/2025"
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = "TARIFA - BB"

$ws.Range("A5").Value = "GRUPO XYZ"
$ws.Range("B5").Value = "BB"
$ws.Range("C5").Value = 1234
$ws.Range("D5").Value = "4567-8"
$ws.Range("E5").Value = "09/03/2025"
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = "TARIFA - BB"

$ws.Range("A6").Value = "GRUPO XYZ"
$ws.Range("B6").Value = "BB"
$ws.Range("C6").Value = 1234
$ws.Range("D6").Value = "4567-8"
$ws.Range("E6").Value = "09/04/2025"
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = "TARIFA - BB"

# --- BRADESCO bank-fee entries ---
$ws.Range("A7").Value = "GRUPO XYZ"
$ws.Range("B7").Value = "BRADESCO"
$ws.Range("C7").Value = 8989
$ws.Range("D7").Value = "8765-4"
$ws.Range("E7").Value = "09/01/2025"
$ws.Range("F7").Value = 13
$ws.Range("G7").Value = "TARIFA - BRADESCO"

$ws.Range("A8").Value = "GRUPO XYZ"
$ws.Range("B8").Value = "BRADESCO"
$ws.Range("C8").Value = 8989
$ws.Range("D8").Value = "8765-4"
$ws.Range("E8").Value = "09/02/2025"
$ws.Range("F8").Value = 14
$ws.Range("G8").Value = "TARIFA - BRADESCO"

# Reflect where the author was last working: "Tarifas bancárias" tab active,
# cursor resting on D2.
$ws.Activate()
$ws.Range("D2").Select()
